$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "rxxx"
$ws.Range("B3").Value = "cameron"
$ws.Range("C3").Value = "it didnt work at all"
$ws.Range("D3").Value = "2025-09-27 00:40:02"
